$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to hold a plain text value (matches the source file's
# use of inline/shared strings for the "Price" column, even when the text
# happens to look like a number, e.g. "250.34"). A direct .Value assignment
# would let Excel auto-coerce such strings into numeric cells, and forcing
# text via NumberFormat="@" leaves a stray (if unused) style behind. Using a
# literal-text formula and then collapsing it to a static value via
# Copy/PasteSpecial(xlPasteValues) yields a plain text cell with no format
# changes at all.
function Set-TextValue($row, $col, $val) {
    $escaped = $val -replace '"', '""'
    $ws.Cells.Item($row, $col).Formula = '="' + $escaped + '"'
    $ws.Cells.Item($row, $col).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial(-4163)
}

# Row => D(new), E(new) values. $null for D means that column is unchanged
# for that row in this edit.
$updates = @(
    @{ Row = 2;  D = "42.806.70"; E = "  -1.09%  " },
    @{ Row = 3;  D = "2.275.66";  E = "  -0.96%  " },
    @{ Row = 4;  D = $null;       E = "  -0.16%  " },
    @{ Row = 5;  D = "250.34";    E = "  -0.96%  " },
    @{ Row = 6;  D = "0.640";     E = "  -0.43%  " },
    @{ Row = 7;  D = "79.11";     E = "  +7.73%  " },
    @{ Row = 8;  D = $null;       E = "  +0.01%  " },
    @{ Row = 9;  D = "0.648";     E = "  -2.64%  " },
    @{ Row = 10; D = "41.47";     E = "  +4.17%  " },
    @{ Row = 11; D = "0.0969";    E = "  -1.45%  " },
    @{ Row = 12; D = "7.39";      E = "  -2.19%  " },
    @{ Row = 13; D = "0.105";     E = "  -0.75%  " },
    @{ Row = 14; D = "2.615.60";  E = "  -0.89%  " },
    @{ Row = 15; D = "15.11";     E = "  -0.68%  " },
    @{ Row = 16; D = "0.873";     E = "  -3.18%  " },
    @{ Row = 17; D = "2.275.95";  E = "  -0.56%  " },
    @{ Row = 18; D = "42.696.64"; E = "  -1.13%  " },
    @{ Row = 19; D = $null;       E = "  -1.68%  " },
    @{ Row = 20; D = $null;       E = "  -3.24%  " },
    @{ Row = 21; D = "72.05";     E = "  -2.41%  " },
    @{ Row = 22; D = "232.20";    E = "  -2.72%  " },
    @{ Row = 23; D = $null;       E = "  -1.30%  " },
    @{ Row = 24; D = $null;       E = "  -3.16%  " },
    @{ Row = 25; D = $null;       E = "  -0.08%  " },
    @{ Row = 26; D = "11.41";     E = "  -4.13%  " },
    @{ Row = 27; D = $null;       E = "  -4.94%  " },
    @{ Row = 28; D = $null;       E = "  +1.35%  " },
    @{ Row = 29; D = "169.66";    E = "  +0.78%  " },
    @{ Row = 30; D = "6.75";      E = "  +6.77%  " },
    @{ Row = 31; D = $null;       E = "  -2.28%  " },
    @{ Row = 32; D = "0.0854";    E = "  +5.17%  " },
    @{ Row = 33; D = $null;       E = "  -4.70%  " },
    @{ Row = 34; D = "30.61";     E = "  -2.34%  " },
    @{ Row = 35; D = $null;       E = "  +0.12%  " },
    @{ Row = 36; D = "4.59";      E = "  -4.23%  " },
    @{ Row = 37; D = $null;       E = "  -1.14%  " },
    @{ Row = 39; D = "13.53";     E = "  -0.47%  " },
    @{ Row = 40; D = $null;       E = "  -3.83%  " },
    @{ Row = 41; D = "5.98";      E = "  -2.58%  " },
    @{ Row = 42; D = "115.88";    E = "  +18.28%  " },
    @{ Row = 43; D = $null;       E = "  -2.20%  " },
    @{ Row = 44; D = "61.51";     E = "  -0.87%  " },
    @{ Row = 45; D = "8.89";      E = "  -3.91%  " },
    @{ Row = 46; D = $null;       E = "  -2.08%  " },
    @{ Row = 47; D = "4.59";      E = "  -7.53%  " },
    @{ Row = 48; D = $null;       E = "  -0.18%  " },
    @{ Row = 49; D = $null;       E = "  -4.02%  " },
    @{ Row = 50; D = "1.17";      E = "  -2.53%  " },
    @{ Row = 51; D = $null;       E = "  -2.92%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        Set-TextValue $u.Row 4 $u.D
    }
    # Column E values (percentages with padding spaces) are never valid
    # numbers, so a plain assignment already keeps them as text.
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
